$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.002.29"
$ws.Range("E2").Value = "  +3.17%  "
$ws.Range("D3").Value = "2.627.02"
$ws.Range("E4").Value = "  -0.10%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "595.51"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.61%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "155.15"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +5.17%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("E9").Value = "  +8.69%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.401"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +5.65%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "5.77"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("E13").Value = "  +6.76%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.0000186"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +22.12%  "
$ws.Range("D15").Value = "3.098.75"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "64.876.68"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "2.614.55"
$ws.Range("E17").Value = "  +2.01%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "12.49"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("E19").Value = "  +3.26%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "351.35"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +2.59%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "7.34"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +8.70%  "
$ws.Range("E22").Value = "  +0.22%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "68.23"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("E24").Value = "  +4.93%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "1.64"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +1.96%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "8.05"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").Value = "0.0₃0951"
$ws.Range("E29").Value = "  +12.71%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "523.71"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -4.91%  "
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("E33").Value = "  +2.31%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "5.54"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +8.15%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "6.29"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +6.63%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.425"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +3.70%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "163.88"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("E38").Value = "  +4.92%  "
$ws.Range("E39").Value = "  +5.90%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.998"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  +0.01%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "164.77"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  +3.99%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.0614"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +5.49%  "
$ws.Range("E46").Value = "  +2.74%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "2.20"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +9.35%  "
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("E49").Value = "  +3.45%  "
$ws.Range("E50").Value = "  +2.05%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "19.38"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +3.05%  "
